# Update the "Förändrad" (Changed) date column (C) for data rows 2-12
# from 2023-10-13 (45212) to 2023-10-22 (45221).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
